# "For PCBWAY March 03"
# Update a few BOM rows:
#  - Row 3 (0.1uF capacitor): clear "Manufacturer Lifecycle 1" (J3)
#  - Row 8 (CRCW0603120KFKEA resistor, R1/R2): clear "Revision ID" (E8)
#    and change "Revision State" (F8) to "Unknown server"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear J3 but keep it as an (empty) text cell, like the other empty text
# cells already used in this column throughout the sheet.
$ws.Range("J3").Value = "'"

# E8 previously held a "quoted text" revision id; reset it to a plain
# numeric value first (which drops the quote-prefix text formatting) and
# then clear it, leaving a plain empty cell.
$ws.Range("E8").Value = 5
$ws.Range("E8").ClearContents()

# F8 keeps its quote-prefixed text formatting, just with new text.
$ws.Range("F8").Value = "'Unknown server"
